$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, shifting existing rows 7-15 down to 8-16.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new record, copying the repeated
# attributes from the (now shifted) row below (column by column, using
# Value2 since those columns are plain text/numbers) and then updating
# the values that actually differ (date, volume, prices).
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(7, $col).Value = $ws.Cells.Item(8, $col).Value2
}

$ws.Range("D7").Value = 44679
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 5500
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = 5500
$ws.Range("P7").Value = 5500
